$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'62.556.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.441.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.74%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'574.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.64%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.76%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.67%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.437.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.93%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.56%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.91%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.882.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.87%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.359.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.436.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.84%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'328.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.45%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.50%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.42%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'65.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'635.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0967"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.989"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.97%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.78%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.17%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E37").Value = "'  -1.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.376"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.75%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'18.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.89%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.92%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'146.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'42.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.60%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'2.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.21%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'145.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.49%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0525"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.22%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.73%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'19.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0230"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.60%  "
$ws.Range("E51").Style = "Normal"
